$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "customer_locations"

# 2. Build the new "Total Number of Customers per continent" side table in D1:E4
$ws.Range("E1").Value = "Total Number of Customers"

$ws.Range("D2").Value = "Europe"
$ws.Range("E2").Value = 54

$ws.Range("D3").Value = "North America"
$ws.Range("E3").Value = 21

$ws.Range("D4").Value = "South America"
$ws.Range("E4").Value = 16

# 3. Bold the header cell and the continent labels (matches the existing bold header style)
$ws.Range("E1").Font.Bold = $true
$ws.Range("D2:D4").Font.Bold = $true

# 4. Put a thin box border around the original table (A1:B22) and the new table (D1:E4)
$ws.Range("A1:B22").Borders.LineStyle = 1
$ws.Range("D1:E4").Borders.LineStyle = 1

# 5. Size the new columns similarly to the existing bestfit column B
$ws.Columns("D").ColumnWidth = 12.33
$ws.Columns("E").ColumnWidth = 23.33
